# Update test-case step/expected-result text for TC1, TC2 and TC3.
# The workbook stores these as shared strings, each referenced by three
# rows (10, 20, 30) in columns B (Steps) and D (Expected Results).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStep = "Administrador acessa a funcionalidade de 'Catalogo (Perfis) de Competencias' a partir do menu inicial"
$newExpected = "SYSTEM exibe a listagem do Catalogo (Perfis) de Competencias cadastradas com a opcao 'Alterar Gerente' dentre as varias listadas"

foreach ($row in 10, 20, 30) {
    $ws.Range("B$row").Value = $newStep
    $ws.Range("D$row").Value = $newExpected
}
